$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old "CASHFLOW_AUG_*" / "CASHFLOW_DEC_*" rows (13-20); we will
# rewrite this block as a full Aug-Dec 2025 cashflow projection table (rows 13-42).
$ws.Range("A13:J20").ClearContents()

$data = @(
    @('CASHFLOW_AUGUST_VISIT_AVG', 'August 2023/2024 Visit Average', 'Remaining Year Cashflow Projections', 'Historical', 840, $null, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_AUGUST_VISIT_PROJ', 'August 2025 Visit Projection', 'Remaining Year Cashflow Projections', 'Forecast', $null, 1056, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_AUGUST_EXPENSES', 'August 2025 Expected Expenses', 'Remaining Year Cashflow Projections', 'Forecast', $null, 102064, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_AUGUST_REVENUE', 'August 2025 Expected Revenue', 'Remaining Year Cashflow Projections', 'Forecast', $null, 120829, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_AUGUST_PROFIT', 'August 2025 Expected Profit', 'Remaining Year Cashflow Projections', 'Forecast', $null, 18765, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_AUGUST_CASH_POSITION', 'August 2025 End of Month Cash Position', 'Remaining Year Cashflow Projections', 'Forecast', $null, 293765, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_SEPTEMBER_VISIT_AVG', 'September 2023/2024 Visit Average', 'Remaining Year Cashflow Projections', 'Historical', 865, $null, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_SEPTEMBER_VISIT_PROJ', 'September 2025 Visit Projection', 'Remaining Year Cashflow Projections', 'Forecast', $null, 1087, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_SEPTEMBER_EXPENSES', 'September 2025 Expected Expenses', 'Remaining Year Cashflow Projections', 'Forecast', $null, 105102, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_SEPTEMBER_REVENUE', 'September 2025 Expected Revenue', 'Remaining Year Cashflow Projections', 'Forecast', $null, 124425, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_SEPTEMBER_PROFIT', 'September 2025 Expected Profit', 'Remaining Year Cashflow Projections', 'Forecast', $null, 19324, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_SEPTEMBER_CASH_POSITION', 'September 2025 End of Month Cash Position', 'Remaining Year Cashflow Projections', 'Forecast', $null, 313089, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_OCTOBER_VISIT_AVG', 'October 2023/2024 Visit Average', 'Remaining Year Cashflow Projections', 'Historical', 915, $null, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_OCTOBER_VISIT_PROJ', 'October 2025 Visit Projection', 'Remaining Year Cashflow Projections', 'Forecast', $null, 1150, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_OCTOBER_EXPENSES', 'October 2025 Expected Expenses', 'Remaining Year Cashflow Projections', 'Forecast', $null, 111177, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_OCTOBER_REVENUE', 'October 2025 Expected Revenue', 'Remaining Year Cashflow Projections', 'Forecast', $null, 131618, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_OCTOBER_PROFIT', 'October 2025 Expected Profit', 'Remaining Year Cashflow Projections', 'Forecast', $null, 20441, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_OCTOBER_CASH_POSITION', 'October 2025 End of Month Cash Position', 'Remaining Year Cashflow Projections', 'Forecast', $null, 333530, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_NOVEMBER_VISIT_AVG', 'November 2023/2024 Visit Average', 'Remaining Year Cashflow Projections', 'Historical', 950, $null, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_NOVEMBER_VISIT_PROJ', 'November 2025 Visit Projection', 'Remaining Year Cashflow Projections', 'Forecast', $null, 1194, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_NOVEMBER_EXPENSES', 'November 2025 Expected Expenses', 'Remaining Year Cashflow Projections', 'Forecast', $null, 115430, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_NOVEMBER_REVENUE', 'November 2025 Expected Revenue', 'Remaining Year Cashflow Projections', 'Forecast', $null, 136652, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_NOVEMBER_PROFIT', 'November 2025 Expected Profit', 'Remaining Year Cashflow Projections', 'Forecast', $null, 21223, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_NOVEMBER_CASH_POSITION', 'November 2025 End of Month Cash Position', 'Remaining Year Cashflow Projections', 'Forecast', $null, 354752, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_DECEMBER_VISIT_AVG', 'December 2023/2024 Visit Average', 'Remaining Year Cashflow Projections', 'Historical', 1175, $null, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_DECEMBER_VISIT_PROJ', 'December 2025 Visit Projection', 'Remaining Year Cashflow Projections', 'Forecast', $null, 1477, 'NULL', 'NULL', 'Count', 'Dashboard Template Cashflow'),
    @('CASHFLOW_DECEMBER_EXPENSES', 'December 2025 Expected Expenses', 'Remaining Year Cashflow Projections', 'Forecast', $null, 142768, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_DECEMBER_REVENUE', 'December 2025 Expected Revenue', 'Remaining Year Cashflow Projections', 'Forecast', $null, 169017, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_DECEMBER_PROFIT', 'December 2025 Expected Profit', 'Remaining Year Cashflow Projections', 'Forecast', $null, 26249, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
    @('CASHFLOW_DECEMBER_CASH_POSITION', 'December 2025 End of Month Cash Position', 'Remaining Year Cashflow Projections', 'Forecast', $null, 381001, 'NULL', 'NULL', 'Currency', 'Dashboard Template Cashflow'),
)

$startRow = 13
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $val = $rowValues[$c]
        if ($val -ne $null) {
            $ws.Cells.Item($r, $c + 1).Value = $val
        }
    }
}

# Reflect the new used range and cursor position recorded in the saved file.
$ws.Range("B7").Select()
